$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain text so numeric-looking values
# (e.g. "17.44", "0.9997", "29.349.50") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.349.50"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.841.22"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "239.87"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "0.9997"
$ws.Range("D8").Value = "0.07457"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "0.2896"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "24.91"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("D11").Value = "0.07727"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "1.843.48"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "4.974"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "0.00001035"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "81.81"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "29.324.80"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "229.06"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "12.33"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "7.373"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "158.03"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "8.524"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "0.1349"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").Value = "17.44"
$ws.Range("D28").Value = "0.06915"
$ws.Range("E28").Value = "  +6.57%  "
$ws.Range("D29").Value = "1.453"
$ws.Range("E29").Value = "  +4.04%  "
$ws.Range("D30").Value = "1.486"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "4.066"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "4.057"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "1.831"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "0.6972"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").Value = "2.582"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "2.821"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").Value = "1.237.33"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("D40").Value = "6.786"
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "101.08"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.974.87"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "65.36"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("E46").Value = "  +4.60%  "
$ws.Range("D47").Value = "7.034"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "1.712"
$ws.Range("E48").Value = "  +2.50%  "
$ws.Range("D49").Value = "8.981"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "0.1142"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").Value = "0.3909"
$ws.Range("E51").Value = "  -0.87%  "
